$wb = $excel.ActiveWorkbook

# ======================================================================
# Part A - "总计" (summary) sheet: push the existing 2021-Q4 summary row
# down to row 3, then turn row 2 into the new 2022-Q3 summary row.
# ======================================================================
$totalSheet = $wb.Worksheets.Item(1)

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 5
$totalSheet.Range("D3").Value = 0.48

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.63

# ======================================================================
# Part B - Rename the current "2021-Q4" sheet to "2022-Q3" (it keeps its
# sheetId/tab slot) and add a brand-new sheet right after it, named
# "2021-Q4", which will receive the untouched old data. This reproduces
# the target sheetId ordering: 总计=1, 2022-Q3=2, 2021-Q4=3.
# ======================================================================
$oldQ4 = $wb.Worksheets.Item("2021-Q4")
$oldQ4.Name = "2022-Q3"

$newQ4 = $wb.Worksheets.Add($null, $oldQ4)
$newQ4.Name = "2021-Q4"
$newQ4.PageSetup.LeftMargin = 54
$newQ4.PageSetup.RightMargin = 54
$newQ4.PageSetup.TopMargin = 72
$newQ4.PageSetup.BottomMargin = 72
$newQ4.PageSetup.HeaderMargin = 36
$newQ4.PageSetup.FooterMargin = 36

# ======================================================================
# Part C - Populate the new "2021-Q4" sheet with the fund-holding data
# that used to live on the old sheet (unchanged content), copying the
# existing header / index-column formatting over first.
# ======================================================================
$oldQ4.Range("B1:H1").Copy()
$newQ4.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$oldQ4.Range("A2").Copy()
$newQ4.Range("A2:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newQ4.Range("B1").Value = "基金代码"
$newQ4.Range("C1").Value = "基金名称"
$newQ4.Range("D1").Value = "基金规模"
$newQ4.Range("E1").Value = "股票总仓位"
$newQ4.Range("F1").Value = "仓位占比"
$newQ4.Range("G1").Value = "持有市值(亿元)"
$newQ4.Range("H1").Value = "仓位排名"

$newQ4.Range("B2:G5").NumberFormat = "@"
$newQ4.Range("B6:F6").NumberFormat = "@"
# (G6 is left as General/numeric - see below, it holds the number 0)

$newQ4.Range("A2").Value = 0
$newQ4.Range("B2").Value = "008526"
$newQ4.Range("C2").Value = "华泰柏瑞行业精选混合A"
$newQ4.Range("D2").Value = "8.59"
$newQ4.Range("E2").Value = "86.59"
$newQ4.Range("F2").Value = "4.24"
$newQ4.Range("G2").Value = "0.3642"
$newQ4.Range("H2").Value = 8

$newQ4.Range("A3").Value = 1
$newQ4.Range("B3").Value = "001601"
$newQ4.Range("C3").Value = "鑫元鑫新收益灵活配置混合A"
$newQ4.Range("D3").Value = "2.01"
$newQ4.Range("E3").Value = "85.37"
$newQ4.Range("F3").Value = "3.81"
$newQ4.Range("G3").Value = "0.0766"
$newQ4.Range("H3").Value = 5

$newQ4.Range("A4").Value = 2
$newQ4.Range("B4").Value = "001537"
$newQ4.Range("C4").Value = "中加改革红利灵活配置混合"
$newQ4.Range("D4").Value = "1.03"
$newQ4.Range("E4").Value = "92.75"
$newQ4.Range("F4").Value = "3.45"
$newQ4.Range("G4").Value = "0.0355"
$newQ4.Range("H4").Value = 5

$newQ4.Range("A5").Value = 3
$newQ4.Range("B5").Value = "008527"
$newQ4.Range("C5").Value = "华泰柏瑞行业精选混合C"
$newQ4.Range("D5").Value = "0.18"
$newQ4.Range("E5").Value = "86.59"
$newQ4.Range("F5").Value = "4.24"
$newQ4.Range("G5").Value = "0.0076"
$newQ4.Range("H5").Value = 8

$newQ4.Range("A6").Value = 4
$newQ4.Range("B6").Value = "001602"
$newQ4.Range("C6").Value = "鑫元鑫新收益灵活配置混合C"
$newQ4.Range("D6").Value = "0.00"
$newQ4.Range("E6").Value = "85.37"
$newQ4.Range("F6").Value = "3.81"
$newQ4.Range("G6").Value = 0
$newQ4.Range("H6").Value = 5

# ======================================================================
# Part D - Overwrite the (renamed) "2022-Q3" sheet with the new fund
# holding data, removing the now-unused 6th row left over from the old
# content.
# ======================================================================
$oldQ4.Rows(6).Delete()

$oldQ4.Range("B1").Value = "基金代码"
$oldQ4.Range("C1").Value = "基金名称"
$oldQ4.Range("D1").Value = "基金规模"
$oldQ4.Range("E1").Value = "股票总仓位"
$oldQ4.Range("F1").Value = "仓位占比"
$oldQ4.Range("G1").Value = "持有市值(亿元)"
$oldQ4.Range("H1").Value = "仓位排名"

$oldQ4.Range("B2:G4").NumberFormat = "@"
$oldQ4.Range("B5:F5").NumberFormat = "@"
# (G5 is left as General/numeric - it holds the number 0)

$oldQ4.Range("A2").Value = 0
$oldQ4.Range("B2").Value = "470009"
$oldQ4.Range("C2").Value = "汇添富民营活力混合A"
$oldQ4.Range("D2").Value = "24.14"
$oldQ4.Range("E2").Value = "90.04"
$oldQ4.Range("F2").Value = "2.53"
$oldQ4.Range("G2").Value = "0.6107"
$oldQ4.Range("H2").Value = 10

$oldQ4.Range("A3").Value = 1
$oldQ4.Range("B3").Value = "014831"
$oldQ4.Range("C3").Value = "兴银中证1000指数增强A"
$oldQ4.Range("D3").Value = "1.37"
$oldQ4.Range("E3").Value = "83.33"
$oldQ4.Range("F3").Value = "0.81"
$oldQ4.Range("G3").Value = "0.0111"
$oldQ4.Range("H3").Value = 7

$oldQ4.Range("A4").Value = 2
$oldQ4.Range("B4").Value = "014832"
$oldQ4.Range("C4").Value = "兴银中证1000指数增强C"
$oldQ4.Range("D4").Value = "0.90"
$oldQ4.Range("E4").Value = "83.33"
$oldQ4.Range("F4").Value = "0.81"
$oldQ4.Range("G4").Value = "0.0073"
$oldQ4.Range("H4").Value = 7

$oldQ4.Range("A5").Value = 3
$oldQ4.Range("B5").Value = "960014"
$oldQ4.Range("C5").Value = "汇添富民营活力混合 O"
$oldQ4.Range("D5").Value = "0.00"
$oldQ4.Range("E5").Value = "90.04"
$oldQ4.Range("F5").Value = "2.53"
$oldQ4.Range("G5").Value = 0
$oldQ4.Range("H5").Value = 10

# Leave the workbook's active selection on the first sheet, matching the
# original file's state.
$totalSheet.Activate() | Out-Null
$totalSheet.Range("A1").Select() | Out-Null
